$wb = $excel.ActiveWorkbook

# --- 1. Add a "State" column to the hotel_info sheet, between Hotel_Name and City ---
$hotel = $wb.Worksheets.Item("hotel_info")

# Column C currently holds "City" / "New Orleans"; insert a new blank column there,
# shifting City (and everything after it) one column to the right.
$hotel.Columns.Item(3).Insert()

$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# --- 2. Reorder the sheet tabs: review_info should come before hotel_info ---
$review = $wb.Worksheets.Item("review_info")
$review.Move($wb.Worksheets.Item(1))

# Keep the first tab active (matches the saved workbookView's activeTab=0)
$wb.Worksheets.Item(1).Activate()
